$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.290976762771606
$ws.Range("B1").Value = 2.369591236114502
$ws.Range("C1").Value = 3.029009103775024
$ws.Range("D1").Value = 3.563786745071411
$ws.Range("E1").Value = 1.201399803161621
